# MitsosBarton2006Ex313 - "alpha_zero" generator run update
# "expermits todos no convexos menos el 5to"
#
# Updates the numeric leader/follower constraint data (x/y shifted from 1 to
# 2.3 / 4.45 respectively) and the derived vectors (vec_bf / vec_BF) that
# depend on them, across the worksheets that hold this experiment's data.

$wb = $excel.ActiveWorkbook

# NOTE: "Vector_bf" and "Vector_BF" differ only by letter case, and sheet
# names are resolved case-insensitively by Worksheets.Item(<name>), so all
# sheets below are addressed by their (1-based) position to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---------------------------------------------------------------------
# Restricciones_del_lider (leader restrictions)
# ---------------------------------------------------------------------
$wsLider = $wb.Worksheets.Item(2)

$wsLider.Range("A2").Value = "2.3000000000000003 - x"
$wsLider.Range("B2").NumberFormat = "@"
$wsLider.Range("B2").Value = "-3.3000000000000003"
$wsLider.Range("D2").NumberFormat = "@"
$wsLider.Range("D2").Value = "0.51"

$wsLider.Range("A3").Value = "-2.3000000000000003 + x"
$wsLider.Range("B3").NumberFormat = "@"
$wsLider.Range("B3").Value = "1.3000000000000003"
$wsLider.Range("D3").NumberFormat = "@"
$wsLider.Range("D3").Value = "0.17"

# ---------------------------------------------------------------------
# Restricciones_del_follower (follower restrictions)
# ---------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item(3)

$wsFollower.Range("A2").Value = "-4.449999999999999 + y"
$wsFollower.Range("B2").NumberFormat = "@"
$wsFollower.Range("B2").Value = "3.4499999999999993"
$wsFollower.Range("D2").NumberFormat = "@"
$wsFollower.Range("D2").Value = "0.82"
$wsFollower.Range("E2").NumberFormat = "@"
$wsFollower.Range("E2").Value = "6.7"
$wsFollower.Range("F2").NumberFormat = "@"
$wsFollower.Range("F2").Value = "1.7000000000000002"

$wsFollower.Range("A3").Value = "4.449999999999999 - y"
$wsFollower.Range("B3").NumberFormat = "@"
$wsFollower.Range("B3").Value = "-5.449999999999999"
$wsFollower.Range("D3").NumberFormat = "@"
$wsFollower.Range("D3").Value = "0.81"
$wsFollower.Range("E3").NumberFormat = "@"
$wsFollower.Range("E3").Value = "8.4"
# F3 ("J_0_LP_v" Gamma_value) stays "0" - unchanged by this edit.

# ---------------------------------------------------------------------
# Punto_modificado (modified point x / y)
# ---------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)

$wsPunto.Range("A2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "2.3000000000000003"
$wsPunto.Range("B2").NumberFormat = "@"
$wsPunto.Range("B2").Value = "4.449999999999999"

# ---------------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------------
$wsVecbf = $wb.Worksheets.Item(5)

$wsVecbf.Range("A2").NumberFormat = "@"
$wsVecbf.Range("A2").Value = "1.9220000000000041"

# ---------------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------------
$wsVecBF = $wb.Worksheets.Item(6)

$wsVecBF.Range("A2").NumberFormat = "@"
$wsVecBF.Range("A2").Value = "-0.66"
$wsVecBF.Range("A3").NumberFormat = "@"
$wsVecBF.Range("A3").Value = "2.7"
